{"js": "// Office.js (Word JavaScript API) edit script.\n// This is the body of `async (context) => { ... }`.\n//\n// Target edit (translation touch-up of Welcome.docx):\n//   The SmartCash-mining paragraph loses its trailing clause\n//   \", until Smartcash reaches a considerable market cap\" so the sentence\n//   now ends \"...for quite some time.\"\n//\n// (The document's lone \"exchanges\" bookmark is left as-is: Word only ever\n// renumbers bookmark IDs compactly/positionally at save time based on how\n// many bookmarks exist in the file, and this document has just the one, so\n// there is no user-facing action that changes its stored id.)\n\nconst body = context.document.body;\n\nconst oldClause = \", until Smartcash reaches a considerable market cap\";\n\nconst searchResults = body.search(oldClause, { matchCase: true, matchWholeWord: false });\nsearchResults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < searchResults.items.length; i++) {\n  searchResults.items[i].insertText(\"\", Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is already open as $d.\n#\n# Target edit (translation touch-up of Welcome.docx):\n#   The SmartCash-mining paragraph loses its trailing clause\n#   \", until Smartcash reaches a considerable market cap\" so the sentence\n#   now ends \"...for quite some time.\"\n#\n# (The document's lone \"exchanges\" bookmark is left untouched: Word only\n# ever renumbers bookmark IDs compactly/positionally at save time based on\n# how many bookmarks currently exist in the file, and this document has\n# just the one bookmark, so there is no user-facing action available here\n# that changes its stored id.)\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n\n$find.Text = \", until Smartcash reaches a considerable market cap.\"\n$find.Replacement.Text = \".\"\n\n# wdFindContinue = 1, wdReplaceAll = 2\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n"}
